# Update res_bus/vm_pu.xlsx results for "case with 380 kV done":
# - Column B (bus 0 voltage setpoint) changes from 1.05 to 1.02 p.u. for all data rows.
# - Columns C-F and I-N (other bus voltage results) are updated with the newly
#   computed power-flow results for the 380 kV case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: constant voltage setpoint for bus 0, now 1.02 p.u. ---
$ws.Range("B2:B25").Value = 1.02

# --- Columns C,D,E,F,I,J,K,L,M,N: new per-row voltage magnitude results ---
$colLetters = @("C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

$newResults = @(
    @(1.036218748958479, 1.034946298474023, 1.044513338906138, 1.053862728952557, 1.037682781627215, 1.041328262338937, 1.037744313713771, 1.047284177382335, 1.056607525719913, 1.04280706791348),
    @(1.037128146719583, 1.035416389222646, 1.045330185210837, 1.05478162590983, 1.037844123692512, 1.041881735219222, 1.038024600668252, 1.047912262288755, 1.057339280999456, 1.043361326788739),
    @(1.03771719302932, 1.035720933713574, 1.045859665969503, 1.055377292733507, 1.037947589179161, 1.042239839757867, 1.038205604346438, 1.048318949657034, 1.057813222090768, 1.043719939876904),
    @(1.037964971709272, 1.035849049680976, 1.046082479935821, 1.055627967755976, 1.037990861903234, 1.042390378751592, 1.038281610826808, 1.048489985480242, 1.058012572410239, 1.043870692653267),
    @(1.038006583224761, 1.035870565875938, 1.046119904251471, 1.055670072200162, 1.037998114419724, 1.042415654400634, 1.038294367488006, 1.048518706890927, 1.058046050350873, 1.04389600419663),
    @(1.037720503297445, 1.035722645271814, 1.045862642357141, 1.055380641258859, 1.037948168273057, 1.042241851300699, 1.038206620293825, 1.048321234794377, 1.057815885406773, 1.043721954276357),
    @(1.036525958129138, 1.035105091346907, 1.044789203141203, 1.054173050073149, 1.037737500850406, 1.041515316514897, 1.037839112027009, 1.047496384043619, 1.056854732546619, 1.042994387727825),
    @(1.034425707531229, 1.034019750555603, 1.042904831974585, 1.05205345672314, 1.037359156594343, 1.040234893128986, 1.037188801722306, 1.046045048609688, 1.055164540256982, 1.041712145993839),
    @(1.033028768805486, 1.033298223444382, 1.04165349576629, 1.050646096822673, 1.037102178257965, 1.039381224899126, 1.03675350300281, 1.045079015933376, 1.054040170008481, 1.040857265457175),
    @(1.032424660132107, 1.032986298804681, 1.041112837851521, 1.050038067374024, 1.03698978567652, 1.039011578430313, 1.0365646107259, 1.044661090181504, 1.053553897895879, 1.040487094047978),
    @(1.032200384927109, 1.032870513243238, 1.040912192123743, 1.0498124251288, 1.036947870558577, 1.038874275992005, 1.0364943879946, 1.04450591119198, 1.053373364465993, 1.040349596624459),
    @(1.032248487376199, 1.032895346119485, 1.040955223201858, 1.049860816764476, 1.036956869063373, 1.038903727781014, 1.036509453693657, 1.044539195001747, 1.05341208541865, 1.040379090238386),
    @(1.032406119066526, 1.032976726349989, 1.041096248758644, 1.050019411484155, 1.036986324373173, 1.03900022894694, 1.036558807304978, 1.044648261862918, 1.053538973108616, 1.040475728447037),
    @(1.032503256780278, 1.033026877671598, 1.0411831629845, 1.050117154399048, 1.037004450584543, 1.039059686628174, 1.03658920781273, 1.044715469178042, 1.05361716472814, 1.040535270564998),
    @(1.033068877940969, 1.033318935573549, 1.041689402424413, 1.050686478668645, 1.037109613857087, 1.039405757158623, 1.036766030723199, 1.045106760248647, 1.054072454787132, 1.040881832555294),
    @(1.033423885502218, 1.033502271290883, 1.042007269792653, 1.051043967794372, 1.037175280797437, 1.039622838163272, 1.036876839413462, 1.04535230747467, 1.054358204524085, 1.04109922183987),
    @(1.033631030102054, 1.033609256165745, 1.042192790093772, 1.05125261690763, 1.037213475176272, 1.039749457526722, 1.036941433090057, 1.045495566901243, 1.054524934114409, 1.04122602101734),
    @(1.033701673690222, 1.033645743372539, 1.042256066950117, 1.051323783161348, 1.037226480127863, 1.039792631371061, 1.036963451188399, 1.045544420712446, 1.054581794154363, 1.041269256173489),
    @(1.033385788815733, 1.033482596099384, 1.041973153866292, 1.051005598954619, 1.037168246526248, 1.03959954745047, 1.036864954731817, 1.045325958860562, 1.054327540416764, 1.041075898051584),
    @(1.032359697204956, 1.032952759756526, 1.041054715301644, 1.049972703565689, 1.036977655137446, 1.038971811730176, 1.036544275536338, 1.04461614280059, 1.053501605359543, 1.040447270874565),
    @(1.031715234001577, 1.032620077809497, 1.040478291175176, 1.049324479483961, 1.036856853933592, 1.038577134068775, 1.036342306779182, 1.044170185074669, 1.052982826786262, 1.040052032725615),
    @(1.032056811150358, 1.032796395811603, 1.040783765911758, 1.049668001207542, 1.03692098455482, 1.038786359353587, 1.036449406567155, 1.044406563838045, 1.053257791405529, 1.04026155513433),
    @(1.033403002837819, 1.033491486322663, 1.041988569033957, 1.051022935774822, 1.037171425344503, 1.039610071522135, 1.036870325029018, 1.045337864554376, 1.054341396032202, 1.041086437068638),
    @(1.034968109120536, 1.034299986430222, 1.043391128284634, 1.052600424465018, 1.03745780753925, 1.040565927373034, 1.037357237304072, 1.046419990369704, 1.055601074252027, 1.042043650344484)
)

$firstRow = 2
for ($i = 0; $i -lt $newResults.Length; $i++) {
    $rowValues = $newResults[$i]
    $targetRow = $firstRow + $i
    for ($j = 0; $j -lt $colLetters.Length; $j++) {
        $ws.Range($colLetters[$j] + $targetRow).Value = $rowValues[$j]
    }
}
